$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RTECreation")
$ws2 = $wb.Worksheets.Item("SearchRTE")

# --- Sheet1 (RTECreation): fix C2 to be text, then add row 3 ---
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "125960302"
$ws1.Range("C2").Style = "Normal"

$ws1.Range("A3").Value = "RTE Job Creation ForCrud"
$ws1.Range("A3").Style = "Normal"
$ws1.Range("B3").Value = "RT00002194"
$ws1.Range("B3").Style = "Normal"

$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = "125960324"
$ws1.Range("C3").Style = "Normal"

# --- Sheet2 (SearchRTE): replace row 2 data, add row 3 ---
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "125960302"
$ws2.Range("A2").Style = "Normal"

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "32391789"
$ws2.Range("B2").Style = "Normal"

$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "3397027"
$ws2.Range("C2").Style = "Normal"

$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "125960313"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "125960324"
$ws2.Range("A3").Style = "Normal"

$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "32391790"
$ws2.Range("B3").Style = "Normal"

$ws2.Range("C3").NumberFormat = "@"
$ws2.Range("C3").Value = "3397028"
$ws2.Range("C3").Style = "Normal"

$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "125960335"
$ws2.Range("D3").Style = "Normal"

# Selections, matching authored state
$ws1.Range("B10").Select()
$ws2.Range("B7").Select()
